$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-7 down to 4-8
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the Black Rock (2nd exposure) record
$ws.Cells.Item(3, 1).Value = "Black Rock"
$ws.Cells.Item(3, 2).Value = "Smile Buffalo Thai restaurant  305 Beach Road, Black Rock VIC 3193"
$ws.Cells.Item(3, 3).Value = "21/12/20 7:30pm-9:30pm"
$ws.Cells.Item(3, 4).Value = "Case dined in restaurant"
